# Fruta / hortaliza, semanal
# Insert a new weekly date block (2021-09-09 / serial 44448) with three
# quality rows (Extra / Primera / Segunda) ahead of the existing data,
# pushing the previous rows 286-293 down to rows 289-296.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at the top of the existing date block (old row 286).
$ws.Range("A286:R288").Insert()

# ---- Row 286: new date (44448), Calidad = Extra ----
$ws.Cells.Item(286,1).Value = 9
$ws.Cells.Item(286,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(286,3).Value = "Metropolitana"
$ws.Cells.Item(286,4).Value2 = 44448
$ws.Cells.Item(286,5).Value = 13
$ws.Cells.Item(286,6).Value = 100112009
$ws.Cells.Item(286,7).Value = "Acelga"
$ws.Cells.Item(286,8).Value = "Sin especificar"
$ws.Cells.Item(286,9).Value = "Extra"
$ws.Cells.Item(286,10).Value = 16
$ws.Cells.Item(286,11).Value = 12000
$ws.Cells.Item(286,12).Value = 13000
$ws.Cells.Item(286,13).Value = 12500
$ws.Cells.Item(286,14).Value = "$/docena de atados"
$ws.Cells.Item(286,15).Value = "Región Metropolitana"
$ws.Cells.Item(286,16).Value = 4167
$ws.Cells.Item(286,17).Value = 3
$ws.Cells.Item(286,18).Value = "Hortaliza"

# ---- Row 287: new date (44448), Calidad = Primera ----
$ws.Cells.Item(287,1).Value = 9
$ws.Cells.Item(287,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(287,3).Value = "Metropolitana"
$ws.Cells.Item(287,4).Value2 = 44448
$ws.Cells.Item(287,5).Value = 13
$ws.Cells.Item(287,6).Value = 100112009
$ws.Cells.Item(287,7).Value = "Acelga"
$ws.Cells.Item(287,8).Value = "Sin especificar"
$ws.Cells.Item(287,9).Value = "Primera"
$ws.Cells.Item(287,10).Value = 43
$ws.Cells.Item(287,11).Value = 10000
$ws.Cells.Item(287,12).Value = 11000
$ws.Cells.Item(287,13).Value = 10488
$ws.Cells.Item(287,14).Value = "$/docena de atados"
$ws.Cells.Item(287,15).Value = "Región Metropolitana"
$ws.Cells.Item(287,16).Value = 3496
$ws.Cells.Item(287,17).Value = 3
$ws.Cells.Item(287,18).Value = "Hortaliza"

# ---- Row 288: new date (44448), Calidad = Segunda ----
$ws.Cells.Item(288,1).Value = 9
$ws.Cells.Item(288,2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(288,3).Value = "Metropolitana"
$ws.Cells.Item(288,4).Value2 = 44448
$ws.Cells.Item(288,5).Value = 13
$ws.Cells.Item(288,6).Value = 100112009
$ws.Cells.Item(288,7).Value = "Acelga"
$ws.Cells.Item(288,8).Value = "Sin especificar"
$ws.Cells.Item(288,9).Value = "Segunda"
$ws.Cells.Item(288,10).Value = 34
$ws.Cells.Item(288,11).Value = 8000
$ws.Cells.Item(288,12).Value = 9000
$ws.Cells.Item(288,13).Value = 8500
$ws.Cells.Item(288,14).Value = "$/docena de atados"
$ws.Cells.Item(288,15).Value = "Región Metropolitana"
$ws.Cells.Item(288,16).Value = 2833
$ws.Cells.Item(288,17).Value = 3
$ws.Cells.Item(288,18).Value = "Hortaliza"
